# Edit described by the commit:
#   1. Change the table (slide 16, graphicFrame "Google Shape;213;p29") table
#      style from {05394B02-D928-4B00-AF89-767FB32815D6} to
#      {D8D41E90-C217-4AEE-B10F-B9B6BD5497BF}.
#   2. Swap the colour scheme carried by the presentation's theme so that the
#      deck's live theme (currently "Integral") becomes the classic
#      "Office Theme" colour values.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 ------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{D8D41E90-C217-4AEE-B10F-B9B6BD5497BF}")
    }
}

# --- 2. Theme colour swap -------------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> target "Office Theme" values
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    # RGB property takes a COLORREF (0x00BBGGRR), so byte-swap the RRGGBB text.
    $bgr = $hex.Substring(4,2) + $hex.Substring(2,2) + $hex.Substring(0,2)
    $tcs.Colors($i).RGB = [Convert]::ToInt32($bgr, 16)
}
